$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (C1 label changes from "repo_name_to_import" stays same text but now
# is placed via shared string index 3; value unaffected since text identical)
$ws.Range("C1").Value = "repo_name_to_import"

# Update data rows: replace existing azure_namespace / repo_name_to_import / github_username
$ws.Range("B2").Value = "repo-migartion/git-project"
$ws.Range("C2").Value = "app-n-pak"
$ws.Range("D2").Value = "anilsb06"

$ws.Range("B3").Value = "repo-migartion/git-project"
$ws.Range("C3").Value = "ALMAtasks"
$ws.Range("D3").Value = "anilsb06"

# New rows 4-6
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "repo-migartion/git-project"
$ws.Range("C4").Value = "CASAshell"
$ws.Range("D4").Value = "anilsb06"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "repo-migartion/git-project"
$ws.Range("C5").Value = "casa-build-utils"
$ws.Range("D5").Value = "anilsb06"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "repo-migartion/git-project"
$ws.Range("C6").Value = "CASAplotserver"
$ws.Range("D6").Value = "anilsb06"

# Column B content got longer ("repo-migartion/git-project"), so the author
# resized/best-fit the column to show it fully.
$ws.Columns("B").ColumnWidth = 24.2

$ws.Range("D6").Select()
